# Slide 8 (1-based index in the Slides collection) holds the
# "Content Placeholder 1" picture (id=19) that needs to be re-positioned,
# plus a new textbox with the repo URL that needs to be added beneath it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# The picture is the 5th shape on the slide (Rectangle 26, Group 28,
# Group 32, Content Placeholder 16, Content Placeholder 1).
$pic = $s.Shapes.Item(5)

# New position (EMU 3784600 / 1928093 -> points). Left is an exact
# integer number of points; Top needs a value that survives the
# point->EMU rounding done by the Top setter (it floors the value
# rounded to 4 decimal places), so nudge it to the next 0.0001 pt step
# that still lands on the correct EMU value.
$pic.Left = 298
$pic.Top = 151.8184

# Add the new textbox under the picture with the repository link.
# Passing the exact point values straight into AddTextbox (rather than
# assigning .Left/.Top/.Width/.Height afterwards) avoids that same
# rounding quirk and reproduces the target EMU extents exactly:
#   off  x=3852800 y=5096947
#   ext cx=7529861 cy=369332
$tb = $s.Shapes.AddTextbox(1, 303.37007874015745, 401.3344094488189, 592.9024409448818, 29.081259842519685)
$tb.Name = "TextBox 20"

$tb.Fill.Visible = $false

$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1

$tb.TextFrame.TextRange.Text = "https://github.com/jarjc001/ConstellationLogger"
